# Checkpoint edit: bump the cached "datetimeFigureOut" footer field from
# 12/29/23 to 12/30/23 across the slide master and every slide layout, and
# fix the "futher" -> "further" typo on slide 25.

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "12/29/23") {
                $tr.Text = "12/30/23"
            }
        }
    }
}

# Slide master's Date Placeholder.
Update-DateField $p.SlideMaster.Shapes

# Every custom layout's Date Placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateField $layouts.Item($j).Shapes
}

# Fix typo "futher" -> "further" on slide 25 (TextBox 33).
$s25 = $p.Slides.Item(25)
for ($i = 1; $i -le $s25.Shapes.Count; $i++) {
    $sh = $s25.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 33" -and $sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "Then what we can futher?") {
            $tr.Text = "Then what we can further?"
        }
    }
}
